$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7860
$ws.Range("D2").Value = 11130
$ws.Range("E2").Value = 13234
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = -32224
